# Rename the Pokemon labels in the header row (row 1, columns B:V) and in
# the header column (column A, rows 2:22). Both ranges list the same 21
# names in the same order, so a single ordered map drives both updates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
    "Voltorb",
    "Koffing",
    "Trevenant",
    "Pikachu",
    "Machoke",
    "Drapion",
    "Pansear",
    "Meditite",
    "Tyrantrum",
    "Snubbull",
    "Wobbuffet",
    "Shelmet",
    "Reshiram",
    "AbomasnowMega Abomasnow",
    "Liepard",
    "Pelipper",
    "Audino",
    "Aipom",
    "Sceptile",
    "KyuremWhite Kyurem",
    "Dugtrio"
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $value = $names[$i]

    # Row 1 runs across columns B (2) through V (22)
    $ws.Cells.Item(1, $i + 2).Value = $value

    # Column A runs down rows 2 through 22
    $ws.Cells.Item($i + 2, 1).Value = $value
}
